# Actualización 10 de Mayo
# Updates the partial-exam statistics tables on the three sheets
# ("1er Parcial", "2o Parcial", "3er Parcial") with refreshed grade
# counts/percentages (columns E:K -> Aprobados, Reprobados, Por_Apro,
# Por_Repro, Promedio, Blancos, Por_Blan).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "1er Parcial"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("1er Parcial")

$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 74.36
$ws.Range("H4").Value = 25.64
$ws.Range("I4").Value = 7.9
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 25.64

$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 71.43000000000001
$ws.Range("H5").Value = 28.57
$ws.Range("I5").Value = 9.199999999999999
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 28.57

$ws.Range("E6").Value = 22
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 88
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 7.8
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 12

$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 76.19
$ws.Range("H8").Value = 23.81
$ws.Range("I8").Value = 7.8
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 23.81

$ws.Range("E9").Value = 22
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 78.56999999999999
$ws.Range("H9").Value = 21.43
$ws.Range("I9").Value = 6.6
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 3.57

$ws.Range("E10").Value = 21
$ws.Range("F10").Value = 7
$ws.Range("G10").Value = 75
$ws.Range("H10").Value = 25
$ws.Range("I10").Value = 6.3
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 14.29

$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = 72.73
$ws.Range("H11").Value = 27.27
$ws.Range("I11").Value = 6.4
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 13.64

$ws.Range("E16").Value = 39
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 7.7
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0

# ---------------------------------------------------------------
# Sheet "2o Parcial"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("2o Parcial")

$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 80.95
$ws.Range("H2").Value = 19.05
$ws.Range("I2").Value = 7.6
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 19.05

$ws.Range("E3").Value = 17
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 80.95
$ws.Range("H3").Value = 19.05
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 19.05

$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 74.36
$ws.Range("H4").Value = 25.64
$ws.Range("I4").Value = 7.8
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 25.64

$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 71.43000000000001
$ws.Range("H5").Value = 28.57
$ws.Range("I5").Value = 9.300000000000001
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 28.57

$ws.Range("E6").Value = 21
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 84
$ws.Range("H6").Value = 16
$ws.Range("I6").Value = 6.8
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 16

$ws.Range("E7").Value = 24
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 96
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 4

$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 76.19
$ws.Range("H8").Value = 23.81
$ws.Range("I8").Value = 7.1
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 23.81

$ws.Range("E9").Value = 17
$ws.Range("F9").Value = 11
$ws.Range("G9").Value = 60.71
$ws.Range("H9").Value = 39.29
$ws.Range("I9").Value = 7.4
$ws.Range("J9").Value = 11
$ws.Range("K9").Value = 39.29

$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 13
$ws.Range("G10").Value = 53.57
$ws.Range("H10").Value = 46.43
$ws.Range("I10").Value = 6.9
$ws.Range("J10").Value = 13
$ws.Range("K10").Value = 46.43

$ws.Range("E11").Value = 18
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 81.81999999999999
$ws.Range("H11").Value = 18.18
$ws.Range("I11").Value = 6.9
$ws.Range("J11").Value = 4
$ws.Range("K11").Value = 18.18

$ws.Range("E16").Value = 27
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 69.23
$ws.Range("H16").Value = 30.77
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 71.43000000000001
$ws.Range("H17").Value = 28.57
$ws.Range("I17").Value = 7.1
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0

# ---------------------------------------------------------------
# Sheet "3er Parcial"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("3er Parcial")

$ws.Range("I2").Value = 8.6
$ws.Range("I3").Value = 8.699999999999999

$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 74.36
$ws.Range("H4").Value = 25.64
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 25.64

$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 71.43000000000001
$ws.Range("H5").Value = 28.57
$ws.Range("I5").Value = 9.300000000000001
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 28.57

$ws.Range("E6").Value = 22
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 88
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 12

$ws.Range("I7").Value = 7.6

$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 76.19
$ws.Range("H8").Value = 23.81
$ws.Range("I8").Value = 7.7
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 23.81

$ws.Range("E9").Value = 22
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 78.56999999999999
$ws.Range("H9").Value = 21.43
$ws.Range("I9").Value = 6.8
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 3.57

$ws.Range("E10").Value = 21
$ws.Range("F10").Value = 7
$ws.Range("G10").Value = 75
$ws.Range("H10").Value = 25
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = 14.29

$ws.Range("E11").Value = 19
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 86.36
$ws.Range("H11").Value = 13.64
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 13.64

$ws.Range("E16").Value = 27
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 69.23
$ws.Range("H16").Value = 30.77
$ws.Range("I16").Value = 7.3
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0

$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 71.43000000000001
$ws.Range("H17").Value = 28.57
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
